# Case and Fatality Demographics Data Updated (2021-06-11 refresh, data pulled 6.10.21)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Cases by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value = 281
$ws.Range("B3").Value = 1387
$ws.Range("B4").Value = 3834
$ws.Range("B5").Value = 15830
$ws.Range("B6").Value = 17349
$ws.Range("B7").Value = 15225
$ws.Range("B8").Value = 12842
$ws.Range("B9").Value = 4648
$ws.Range("B10").Value = 3147
$ws.Range("B11").Value = 1908
$ws.Range("B12").Value = 1257
$ws.Range("B12").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 27187
$ws.Range("B3").Value = 51585
$ws.Range("B4").Value = 899
$ws.Range("B7").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by RaceEthnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 968
$ws.Range("B3").Value = 13075
$ws.Range("B4").Value = 28531
$ws.Range("B5").Value = 573
$ws.Range("B6").Value = 27822
$ws.Range("B7").Value = 8702

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B5").Value = 262
$ws.Range("B6").Value = 880
$ws.Range("B7").Value = 2548
$ws.Range("B8").Value = 5780
$ws.Range("B9").Value = 4797
$ws.Range("B10").Value = 6172
$ws.Range("B11").Value = 6796
$ws.Range("B12").Value = 6702
$ws.Range("B13").Value = 16816
$ws.Range("B2:B14").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 21297
$ws.Range("B3").Value = 29509
$ws.Range("D13").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Race-Ethnicity  (ends up the active tab, per the diff)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1082
$ws.Range("B3").Value = 5161
$ws.Range("B4").Value = 23593
$ws.Range("B5").Value = 277
$ws.Range("B6").Value = 20671
$ws.Activate()
$ws.Range("G7").Select()
